# Update solver config slightly — adjust the mechanical-separation moisture
# scenario inputs on the Dairy and Beef tabs, and move the active
# selection/tab from Dairy to Beef.

$wb = $excel.ActiveWorkbook

$dairy = $wb.Worksheets.Item("Dairy")
$beef  = $wb.Worksheets.Item("Beef")

# --- Beef sheet updates -----------------------------------------------
# Update the note text BEFORE the Dairy note, so the shared-string table
# ends up with the same slot reuse/append order as the authored edit.
$beef.Range("D3").Value = "68% moisture content of feedstock after mechanical separation results in a balanced system"

$beef.Range("B2").Value = 0.68
$beef.Range("B12").Value = 0.01
$beef.Range("B13").Value = 13.72
$beef.Range("B16").NumberFormat = "0.0%"
$beef.Range("B16").Value = 0.465
$beef.Range("B17").Value = 13.75

# --- Dairy sheet updates ------------------------------------------------
$dairy.Range("D3").Value = "72% moisture content of feedstock after mechanical separation results in a balanced system"

$dairy.Range("B2").Value = 0.72
$dairy.Range("B12").Value = 0.01
$dairy.Range("B13").Value = 15.26
$dairy.Range("B17").Value = 15.3

# --- Move active tab / selection from Dairy to Beef ---------------------
$dairy.Activate()
$dairy.Range("B1:B33").Select()

$beef.Activate()
$beef.Range("E21").Select()
